$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook window position (bookViews/workbookView yWindow 600 -> 1200) ---
try {
    $excel.ActiveWindow.Top = 1200
} catch {}

# --- New rows 76-101: Contact / Editions / Index / Presse / Mail / Results Content sections ---
# Row 76
$ws.Range("A76").Value = "// Contact"
$ws.Range("C76").Value = "Contact"
$ws.Range("C76").WrapText = $true

# Row 77
$ws.Range("A77").Value = "contact_title"
$ws.Range("B77").Value = "Contact"
$ws.Range("B77").WrapText = $true

# Row 78
$ws.Range("A78").Value = "message"
$ws.Range("B78").Value = "Message"
$ws.Range("B78").WrapText = $true
$ws.Range("C78").Value = "Message"
$ws.Range("C78").WrapText = $true

# Row 81
$ws.Range("A81").Value = "// Editions"

# Row 82
$ws.Range("A82").Value = "editions_title"
$v1 = @'
<span class="city">Saint-Tropez</span> Fight Night <span class="year">2017</span>
'@
$ws.Range("B82").Value = $v1
$v2 = @'
<span class="city">Saint-Tropez</span> Fight Night <span class="year">2017</span>
'@
$ws.Range("C82").Value = $v2

# Row 83
$ws.Range("A83").Value = "results"
$ws.Range("B83").Value = "Résultats"
$ws.Range("C83").Value = "Results"

# Row 84
$ws.Range("A84").Value = "all_albums"
$ws.Range("B84").Value = "Tous les albums"
$ws.Range("C84").Value = "All albums"

# Row 85
$ws.Range("A85").Value = "// Index"

# Row 86
$ws.Range("A86").Value = "title_year"
$ws.Range("B86").Value = "Title / year"
$ws.Range("C86").Value = "Title / year"

# Row 87
$ws.Range("A87").Value = "discover"
$ws.Range("B87").Value = "Discover"
$ws.Range("C87").Value = "discover"

# Row 89
$ws.Range("A89").Value = "// Presse"

# Row 90
$ws.Range("A90").Value = "presse"
$ws.Range("B90").Value = "Presse"
$ws.Range("C90").Value = "Press"

# Row 91
$ws.Range("A91").Value = "in_presse"
$ws.Range("B91").Value = "La presse en parle"
$ws.Range("C91").Value = "In the press"

# Row 92
$ws.Range("A92").Value = "// Mail"

# Row 93
$ws.Range("A93").Value = "site mail"
$ws.Range("B93").Value = "info@fightnights.pro"

# Row 95
$ws.Range("A95").Value = "// Results Content"

# Row 96
$ws.Range("A96").Value = "results_content_1"
$v3 = @'
<h1 class="title--translated">Résultats 2013</h1>
                    <h3>CHAMPIONNAT DU MONDE WKN</h3>
                    <p>Kick boxing . 5×2 mn . +100kg<br>
                        Vainqueur : Jerome LE BANNER par KO2
                    </p>
                    <p>
                        Kick boxing . 5×2 mn . -100kg<br>
                        Vainqueur : Wladimir MINEEV aux points
                    </p>
                    <p>
                        Kick boxing . 5×2 mn . -96,600kg<br>
                        Vainqueur : Emmanuel PAYET aux points
                    </p>
                    <h3>MASTERFIGHT</h3>
                    <p>
                        PRO FIGHT . 3×2 mn (1Extra Round) . + 95kg<br>
                        Vainqueur : Anthony REA par KO1
                    </p>
                    <p>
                        PRO FIGHT . 3×2 mn (1Extra Round) . -70kg<br>
                        Vainqueur : Philippe SALMON KO1
                    </p>
                    <p>
                        PRO FIGHT . 3×2 mn (1Extra Round) . -95kg<br>
                        Vainqueur : Zinedine Hameur Lain aux points
                    </p>
                    <h3>SUPERFIGHT</h3>
                    <p>
                        K-1 . 4×2 mn . +100kg<br>
                        Vainqueur : Nicolas WAMBA par KO4
                    </p>
                    <h3>SUPERFIGHT (FÉMININ)</h3>
                    <p>
                        K-1 . 4×2 mn . -58,500kg<br>
                        Vainqueur : Alena Ola aux points
                    </p>
                    <p>
                        K-1 . 4×2 mn . +100kg<br>
                        Vainqueur : Luca PANTO aux points
                    </p>
                    <h3>TOURNOI MASTERFIGHT</h3>
                    <p>
                        Finale PRO FIGHT . 3×2 mn (1Extra Round) . -82kg<br>
                        Vainqueur : Yoan KONGOLO par KO2
                    </p>
                    <p>
                        Demi finale PRO FIGHT . 3×2 mn (1Extra Round) . -82kg<br>
                        Vainqueur : Yoan KONGOLO aux points
                    </p>
                    <p>
                        Demi finale PRO FIGHT . 3×2 mn (1Extra Round) . -82kg<br>
                        Vainqueur : Francky BRUCTER aux points
                    </p>
'@
$ws.Range("B96").Value = $v3
$ws.Range("B96").WrapText = $true
$ws.Rows.Item(96).RowHeight = 72.75

# Row 97
$ws.Range("A97").Value = "results_content_2"
$v4 = @'
<h1 class="title--translated">Résultats 2014</h1>
                    <p>
                        + 110.600 KG<br>
                        Victoire Jérôme Le Banner (Fra) par KO2
                    </p>
                    <p>
                        + 110.600 KG<br>
                        Victoire Freddy Kemayo (Fra) arr arbitre 2e
                    </p>
                    <p>
                        -96.400 KG<br>
                        Victoire Vladimir Mineev (Russie) aux Points
                    </p>
                    <p>
                        -110.600 KG<br>
                        Victoire Nicolas Wamba (Fra) par KO2
                    </p>
                    <p>
                        -63.500 KG<br>
                        Victoire Maneenoi Ekkarit (Tha) par KO1
                    </p>
                    <p>
                        -96.400 KG<br>
                        Victoire Filip Verlinden (Bel) aux points
                    </p>
                    <p>
                        -66.100 KG<br>
                        Victoire Samsamut Kietchongkao (Tha) aux points
                    </p>
                    <p>
                        -80.500 KG<br>
                        Victoire Yohann Kongolo (Sui) aux points
                    </p>
                    <p>
                        -58.200 KG<br>
                        Victoire Taehiran Chommanee (Tha) aux points
                    </p>
'@
$ws.Range("B97").Value = $v4
$ws.Range("B97").WrapText = $true
$ws.Rows.Item(97).RowHeight = 102

# Row 98
$ws.Range("A98").Value = "return"
$ws.Range("B98").Value = "Retour aux editions"
$ws.Range("C98").Value = "Back to the editions"

# Row 99
$ws.Range("A99").Value = "results_content_3"
$v5 = @'
<h1 class="title--translated">Résultats 2015</h1>
                    <h3>KICK-BOXING</h3>
                    <p>
                        75KG<br>
                        Sharos Huyer bat Bakari Tounkara aux points
                    </p>
                    <p>
                        100KG<br>
                        Frank Munoz bat Zinedine Hameur Lain  aux points
                    </p>
                    <p>
                        120KG<br>
                        Fabrice Aurieng bat Yuksel Ayadi aux points
                    </p>
                    <p>
                        95KG<br>
                        Alexey Papin bat Danyo Ilunga aux points
                    </p>
                    <p>
                        93KG<br>
                        Alexander Vezhevatov bat Filip Verlinden aux points
                    </p>
                    <h3>CHAMPIONNAT DU MONDE  – WKN</h3>
                    <p>
                        120 KG<br>
                        Jerome LEBANNER bat Karl ROBERSON aux points
                    </p>
                    <h3>MUAYTHAI</h3>
                    <p>
                        71 KG<br>
                        Yodwicha Por Boonsit (THA) bat aux pts Jimmy Viennot (F)
                    </p>
                    <p>
                        77 KG<br>
                        Yohan LIDON bat Jonathan OLIVEIRA par TKO3
                    </p>
                    <h3>CHAMPIONNAT DU MONDE  – WMC EN 52KG</h3>
                    <p>Lizzie LARGILLIERE bat Petchoydying MOR par abandon à l’appel de R5</p>
'@
$ws.Range("B99").Value = $v5
$ws.Range("B99").WrapText = $true
$ws.Rows.Item(99).RowHeight = 93

# Row 100
$ws.Range("A100").Value = "results_content_4"
$v6 = @'
<h1 class="title--translated">Résultats 2016</h1>
                    <h3>SUPERFIGHT . MUAYTHAI . 3×3 MN . -67KG</h3>
                    <p>
                        Vainqueur : Dylan Salvador par TKO2
                    </p>
                    <h3>500KG TOURNAMENT . DEMI FINALE . 3×3 MN . K-1 RULES . +110KG</h3>
                    <p>
                        Vainqueur : Thomas Vanneste aux points
                    </p>
                    <h3>500KG TOURNAMENT . DEMI FINALE . 3×3 MN . K-1 RULES . +110KG</h3>
                    <p>
                        Vainqueur : Tomas Mozny aux points
                    </p>
                    <h3>SUPERFIGHT . K-1 RULES . 3×3 MN . -96KG</h3>
                    <p>
                        Vainqueur : Mikhael Tiuterev aux points
                    </p>
                    <h3>CHAMPIONNAT DU MONDE WMC . MUAYTHAI . 5×3 MN . -57,200KG</h3>
                    <p>Vainqueur : Taiheran Chomanee aux points</p>
                    <h3>SUPERFIGHT . K-1 RULES . 3×3 MN . -96KG</h3>
                    <p>Vainqueur : Stéphane Susperregui aux points</p>
                    <h3>500KG TOURNAMENT . DEMI FINALE . 3×3 MN . K-1 RULES . +110KG</h3>
                    <p>Vainqueur : Tomas Mozny aux points</p>
                    <h3>CHAMPIONNAT DU MONDE WKN . K-1 RULES . 5×3 MN . -76,600KG</h3>
                    <p>Vainqueur : Yohan Lidon aux points</p>
'@
$ws.Range("B100").Value = $v6
$ws.Range("B100").WrapText = $true
$ws.Rows.Item(100).RowHeight = 90.75

# Row 101
$ws.Range("A101").Value = "results_content_5"
$v7 = @'
<h1 class="title--translated">Résultats 2017</h1>
                    <h3>CHAMPIONNAT DU MONDE WKN . K-1  (79,4 KG)</h3>
                    <p>Yohan LIDON bat Florian KROGER par KO au round 4 (highkick)</p>
                    <h3>SUPERFIGHT . MUAYTHAI</h3>
                    <p>Umar SEMATA bat Evgeny KURAVSKOI par décision</p>
                    <h3>SUPERFIGHT . K-1</h3>
                    <p>Grégory TONY bat Bob SAPP par KO au round 1</p>
                    <p>Mikhail CHALYKH  bat Phillip VERLINDEN par décision</p>
                    <p>Stéphane SUSPERREGUI bat Danyo ILUNGA par décision (unanime)</p>
                    <p>Tomas MOZNY bat Daniel SAM par décision</p>
                    <p>Tomas MOZNY bat Daniel SAM par décision</p>
                    <p>Mallaury KALASHNIKOFF bat Marina SPASIC par décision (unanime)</p>
'@
$ws.Range("B101").Value = $v7
$ws.Range("B101").WrapText = $true
$ws.Rows.Item(101).RowHeight = 104.25

# --- Hyperlink on B93 (mail address) ---
$ws.Hyperlinks.Add($ws.Range("B93"), "mailto:info@fightnights.pro", "", "", "info@fightnights.pro")

Write-Output "done"
